# Actualización automática 2025-10-16 17:30:09
#
# This script updates sales figures for "ROCAFUERTE LOPEZ EVELYN ESTEFANIA"
# (row 16 on both sheets) on the "ILLER LOPEZ ROBERTO FERNANDO" workbook,
# and refreshes the corresponding summary/count cells on row 18.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" --------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# New sales figures for ROCAFUERTE LOPEZ EVELYN ESTEFANIA (row 16)
$wsGrupo.Range("E16").Value = 71.55
$wsGrupo.Range("M16").Value = 944.1900000000001

# Update the "x de 16" completion counters on the totals row (row 18)
$wsGrupo.Range("E18").Value = "1 de 16"
$wsGrupo.Range("M18").Value = "6 de 16"

# --- Sheet "VENTA MENSUAL" ------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# New "octubre" sale figure for ROCAFUERTE LOPEZ EVELYN ESTEFANIA (row 16)
$wsMensual.Range("F16").Value = 1015.74

# Updated "octubre" total (row 18)
$wsMensual.Range("F18").Value = 9165.84
